# Refresh Universalis market-price derived figures (currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) across the per-job leve sheets, as
# produced by the scheduled price-update runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 317.08334
$ws.Range("I11").Value = 317.08334
$ws.Range("K11").Value = 317.08334
$ws.Range("M11").Value = -177.08334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 760
$ws.Range("I39").Value = 760
$ws.Range("K39").Value = 2280
$ws.Range("M39").Value = -1984

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 4869.8
$ws.Range("J125").Value = 7566.6665
$ws.Range("L125").Value = 68099.9985
$ws.Range("N125").Value = -73019.9985

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2601.8215
$ws.Range("I132").Value = 2631.6667
$ws.Range("K132").Value = 7895.000100000001
$ws.Range("M132").Value = -5365.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2970
$ws.Range("I45").Value = 1327.75
$ws.Range("K45").Value = 1327.75
$ws.Range("M45").Value = -950.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4485.6
$ws.Range("I63").Value = 2625
$ws.Range("J63").Value = 5162.1816
$ws.Range("K63").Value = 2625
$ws.Range("L63").Value = 5162.1816
$ws.Range("M63").Value = -1939
$ws.Range("N63").Value = -6534.1816

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4485.6
$ws.Range("I66").Value = 2625
$ws.Range("J66").Value = 5162.1816
$ws.Range("K66").Value = 13125
$ws.Range("L66").Value = 25810.908
$ws.Range("M66").Value = -9693
$ws.Range("N66").Value = -32674.908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2862.3235
$ws.Range("I74").Value = 2370.6667
$ws.Range("J74").Value = 6549.75
$ws.Range("K74").Value = 2370.6667
$ws.Range("L74").Value = 6549.75
$ws.Range("M74").Value = -1496.6667
$ws.Range("N74").Value = -8297.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2862.3235
$ws.Range("I77").Value = 2370.6667
$ws.Range("J77").Value = 6549.75
$ws.Range("K77").Value = 11853.3335
$ws.Range("L77").Value = 32748.75
$ws.Range("M77").Value = -7485.333500000001
$ws.Range("N77").Value = -41484.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2250
$ws.Range("I110").Value = 1000
$ws.Range("K110").Value = 1000
$ws.Range("M110").Value = 1045

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3193.3396
$ws.Range("I132").Value = 2466.8538
$ws.Range("J132").Value = 5675.5
$ws.Range("K132").Value = 7400.5614
$ws.Range("L132").Value = 17026.5
$ws.Range("M132").Value = -4870.5614
$ws.Range("N132").Value = -22086.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2706.889
$ws.Range("I105").Value = 2334.8823
$ws.Range("K105").Value = 2334.8823
$ws.Range("M105").Value = -587.8823000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1474.75
$ws.Range("I16").Value = 1474.75
$ws.Range("K16").Value = 1474.75
$ws.Range("M16").Value = -1187.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2830.5
$ws.Range("I58").Value = 2651.5278
$ws.Range("J58").Value = 3904.3333
$ws.Range("K58").Value = 2651.5278
$ws.Range("L58").Value = 3904.3333
$ws.Range("M58").Value = -2448.5278
$ws.Range("N58").Value = -4310.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1474.75
$ws.Range("I113").Value = 1474.75
$ws.Range("K113").Value = 1474.75
$ws.Range("M113").Value = 695.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3047.7778
$ws.Range("I134").Value = 2803.75
$ws.Range("K134").Value = 8411.25
$ws.Range("M134").Value = -5876.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2830.5
$ws.Range("I136").Value = 2651.5278
$ws.Range("J136").Value = 3904.3333
$ws.Range("K136").Value = 7954.5834
$ws.Range("L136").Value = 11712.9999
$ws.Range("M136").Value = -5404.5834
$ws.Range("N136").Value = -16812.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 100
$ws.Range("K8").Value = 300
$ws.Range("M8").Value = -161

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2478.5642
$ws.Range("I132").Value = 2159.3
$ws.Range("J132").Value = 3542.7778
$ws.Range("K132").Value = 6477.900000000001
$ws.Range("L132").Value = 10628.3334
$ws.Range("M132").Value = -3947.900000000001
$ws.Range("N132").Value = -15688.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3860.6
$ws.Range("I7").Value = 4108
$ws.Range("K7").Value = 4108
$ws.Range("M7").Value = -3996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 37042092
$ws.Range("I40").Value = 66669868
$ws.Range("K40").Value = 66669868
$ws.Range("M40").Value = -66669732

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3238
$ws.Range("I46").Value = 687
$ws.Range("J46").Value = 3748.2
$ws.Range("K46").Value = 687
$ws.Range("L46").Value = 3748.2
$ws.Range("M46").Value = -499
$ws.Range("N46").Value = -4124.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 486.22223
$ws.Range("I55").Value = 458.5
$ws.Range("J55").Value = 541.6667
$ws.Range("K55").Value = 458.5
$ws.Range("L55").Value = 541.6667
$ws.Range("M55").Value = -285.5
$ws.Range("N55").Value = -887.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1915.5264
$ws.Range("I61").Value = 1769.5
$ws.Range("J61").Value = 2077.7778
$ws.Range("K61").Value = 1769.5
$ws.Range("L61").Value = 2077.7778
$ws.Range("M61").Value = -1567.5
$ws.Range("N61").Value = -2481.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3495
$ws.Range("I100").Value = 3495
$ws.Range("K100").Value = 3495
$ws.Range("M100").Value = -2954

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1915.5264
$ws.Range("I113").Value = 1769.5
$ws.Range("J113").Value = 2077.7778
$ws.Range("K113").Value = 1769.5
$ws.Range("L113").Value = 2077.7778
$ws.Range("M113").Value = 400.5
$ws.Range("N113").Value = -6417.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3860.6
$ws.Range("I126").Value = 4108
$ws.Range("K126").Value = 12324
$ws.Range("M126").Value = -9854

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2386.8235
$ws.Range("I136").Value = 2184.7693
$ws.Range("K136").Value = 6554.3079
$ws.Range("M136").Value = -4004.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5505.7144
$ws.Range("I62").Value = 4280
$ws.Range("K62").Value = 4280
$ws.Range("M62").Value = -3656

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5505.7144
$ws.Range("I65").Value = 4280
$ws.Range("K65").Value = 21400
$ws.Range("M65").Value = -18280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4588
$ws.Range("I81").Value = 3691.75
$ws.Range("J81").Value = 5783
$ws.Range("K81").Value = 7383.5
$ws.Range("L81").Value = 11566
$ws.Range("M81").Value = -6322.5
$ws.Range("N81").Value = -13688

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 4588
$ws.Range("I84").Value = 3691.75
$ws.Range("J84").Value = 5783
$ws.Range("K84").Value = 36917.5
$ws.Range("L84").Value = 57830
$ws.Range("M84").Value = -31613.5
$ws.Range("N84").Value = -68438

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1177.7297
$ws.Range("I136").Value = 936.08826
$ws.Range("K136").Value = 2808.26478
$ws.Range("M136").Value = -258.26478
